$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 111114080
$ws.Range("I106").Value = 333335500
$ws.Range("J106").Value = 3366.6667
$ws.Range("K106").Value = 333335500
$ws.Range("L106").Value = 3366.6667
$ws.Range("M106").Value = -333334869
$ws.Range("N106").Value = -4628.6667
$ws.Range("H129").Value = 926.1507
$ws.Range("I129").Value = 956.55554
$ws.Range("J129").Value = 921.875
$ws.Range("K129").Value = 2869.66662
$ws.Range("L129").Value = 2765.625
$ws.Range("M129").Value = 2130.33338
$ws.Range("N129").Value = -12765.625
$ws.Range("H132").Value = 11635323
$ws.Range("I132").Value = 15631978
$ws.Range("J132").Value = 8691.637000000001
$ws.Range("K132").Value = 46895934
$ws.Range("L132").Value = 26074.911
$ws.Range("M132").Value = -46893404
$ws.Range("N132").Value = -31134.911
$ws.Range("H137").Value = 1388
$ws.Range("I137").Value = 1257.6957
$ws.Range("K137").Value = 3773.0871
$ws.Range("M137").Value = -1223.0871
$ws.Range("H138").Value = 1363.0312
$ws.Range("I138").Value = 776.74
$ws.Range("J138").Value = 3456.9285
$ws.Range("K138").Value = 2330.22
$ws.Range("L138").Value = 10370.7855
$ws.Range("M138").Value = 2809.78
$ws.Range("N138").Value = -20650.7855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 764.72
$ws.Range("I32").Value = 764.72
$ws.Range("K32").Value = 764.72
$ws.Range("M32").Value = -477.72
$ws.Range("H61").Value = 1052.0151
$ws.Range("I61").Value = 1006
$ws.Range("J61").Value = 1150.619
$ws.Range("K61").Value = 1006
$ws.Range("L61").Value = 1150.619
$ws.Range("M61").Value = -794
$ws.Range("N61").Value = -1574.619
$ws.Range("H74").Value = 1161.2683
$ws.Range("I74").Value = 1181.8918
$ws.Range("K74").Value = 1181.8918
$ws.Range("M74").Value = -307.8918000000001
$ws.Range("H77").Value = 1161.2683
$ws.Range("I77").Value = 1181.8918
$ws.Range("K77").Value = 5909.459000000001
$ws.Range("M77").Value = -1541.459000000001
$ws.Range("H132").Value = 27150596
$ws.Range("I132").Value = 33334298
$ws.Range("J132").Value = 6538258.5
$ws.Range("K132").Value = 100002894
$ws.Range("L132").Value = 19614775.5
$ws.Range("M132").Value = -100000364
$ws.Range("N132").Value = -19619835.5
$ws.Range("H136").Value = 1052.0151
$ws.Range("I136").Value = 1006
$ws.Range("J136").Value = 1150.619
$ws.Range("K136").Value = 3018
$ws.Range("L136").Value = 3451.857
$ws.Range("M136").Value = -468
$ws.Range("N136").Value = -8551.857

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 52633424
$ws.Range("I105").Value = 1790.6666
$ws.Range("K105").Value = 1790.6666
$ws.Range("M105").Value = -43.66660000000002
$ws.Range("H107").Value = 55556636
$ws.Range("I107").Value = 83334060
$ws.Range("J107").Value = 1771
$ws.Range("K107").Value = 83334060
$ws.Range("L107").Value = 1771
$ws.Range("M107").Value = -83332140
$ws.Range("N107").Value = -5611
$ws.Range("H134").Value = 4635132
$ws.Range("I134").Value = 956.9231
$ws.Range("J134").Value = 10111884
$ws.Range("K134").Value = 2870.7693
$ws.Range("L134").Value = 30335652
$ws.Range("M134").Value = -335.7692999999999
$ws.Range("N134").Value = -30340722

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 12422.223
$ws.Range("I22").Value = 240
$ws.Range("J22").Value = 27650
$ws.Range("K22").Value = 240
$ws.Range("L22").Value = 27650
$ws.Range("M22").Value = 110
$ws.Range("N22").Value = -28350
$ws.Range("H31").Value = 1139.011
$ws.Range("I31").Value = 931.7655999999999
$ws.Range("J31").Value = 1630.2593
$ws.Range("K31").Value = 931.7655999999999
$ws.Range("L31").Value = 1630.2593
$ws.Range("M31").Value = -636.7655999999999
$ws.Range("N31").Value = -2220.2593
$ws.Range("H34").Value = 1139.011
$ws.Range("I34").Value = 931.7655999999999
$ws.Range("J34").Value = 1630.2593
$ws.Range("K34").Value = 931.7655999999999
$ws.Range("L34").Value = 1630.2593
$ws.Range("M34").Value = -729.7655999999999
$ws.Range("N34").Value = -2034.2593
$ws.Range("H100").Value = 35000
$ws.Range("J100").Value = 35000
$ws.Range("L100").Value = 35000
$ws.Range("N100").Value = -37164
$ws.Range("H134").Value = 1273.1632
$ws.Range("I134").Value = 1102.75
$ws.Range("J134").Value = 1593.9412
$ws.Range("K134").Value = 3308.25
$ws.Range("L134").Value = 4781.8236
$ws.Range("M134").Value = -773.25
$ws.Range("N134").Value = -9851.8236

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 223.90909
$ws.Range("I6").Value = 195.375
$ws.Range("J6").Value = 300
$ws.Range("K6").Value = 586.125
$ws.Range("L6").Value = 900
$ws.Range("M6").Value = -473.125
$ws.Range("N6").Value = -1126
$ws.Range("H9").Value = 183367040
$ws.Range("I9").Value = 50000250
$ws.Range("J9").Value = 250050430
$ws.Range("K9").Value = 150000750
$ws.Range("L9").Value = 750151290
$ws.Range("M9").Value = -150000526
$ws.Range("N9").Value = -750151738
$ws.Range("H47").Value = 472
$ws.Range("I47").Value = 268.33334
$ws.Range("J47").Value = 624.75
$ws.Range("K47").Value = 805.0000200000001
$ws.Range("L47").Value = 1874.25
$ws.Range("M47").Value = -374.0000200000001
$ws.Range("N47").Value = -2736.25
$ws.Range("H131").Value = 972.62
$ws.Range("I131").Value = 661.8
$ws.Range("J131").Value = 988.97894
$ws.Range("K131").Value = 1985.4
$ws.Range("L131").Value = 2966.93682
$ws.Range("M131").Value = 3054.6
$ws.Range("N131").Value = -13046.93682
$ws.Range("H137").Value = 41667696
$ws.Range("I137").Value = 29412824
$ws.Range("K137").Value = 88238472
$ws.Range("M137").Value = -88233372

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1560
$ws.Range("I46").Value = 740
$ws.Range("J46").Value = 2995
$ws.Range("K46").Value = 740
$ws.Range("L46").Value = 2995
$ws.Range("M46").Value = -552
$ws.Range("N46").Value = -3371

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 926
$ws.Range("I81").Value = 836.3333
$ws.Range("K81").Value = 1672.6666
$ws.Range("M81").Value = -611.6666
$ws.Range("H84").Value = 926
$ws.Range("I84").Value = 836.3333
$ws.Range("K84").Value = 8363.333000000001
$ws.Range("M84").Value = -3059.333000000001
$ws.Range("H132").Value = 6188245.5
$ws.Range("I132").Value = 20239.963
$ws.Range("J132").Value = 18524256
$ws.Range("K132").Value = 60719.889
$ws.Range("L132").Value = 55572768
$ws.Range("M132").Value = -58189.889
$ws.Range("N132").Value = -55577828
$ws.Range("H136").Value = 7940378.5
$ws.Range("I136").Value = 11909658
$ws.Range("J136").Value = 1819.2858
$ws.Range("K136").Value = 35728974
$ws.Range("L136").Value = 5457.857400000001
$ws.Range("M136").Value = -35726424
$ws.Range("N136").Value = -10557.8574
